$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header strings (volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"


# Row 15
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -80
$ws.Range("L15").Value = -66.666666666666
$ws.Range("M15").Value = 0
$ws.Range("M15").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -30
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = 8.571428571428
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = 8.333333333333
$ws.Range("L16").Value = 23.809523809523
$ws.Range("M16").Value = -3.703703703703
$ws.Range("N16").Value = -76.991150442477

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 18
$ws.Range("J17").Value = 32
$ws.Range("K17").Value = -43.75
$ws.Range("L17").Value = -40
$ws.Range("M17").Value = -30.769230769230
$ws.Range("N17").Value = -57.142857142857

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -23.809523809523
$ws.Range("I18").Value = 15
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = 15.384615384615
$ws.Range("L18").Value = 15.384615384615
$ws.Range("M18").Value = -58.333333333333
$ws.Range("N18").Value = -88.888888888888

# Row 19
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -42.105263157894
$ws.Range("G19").Value = 73
$ws.Range("H19").Value = -32.876712328767
$ws.Range("I19").Value = 32
$ws.Range("J19").Value = 44
$ws.Range("K19").Value = -27.272727272727
$ws.Range("L19").Value = -17.948717948717
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -46.666666666666

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 30
$ws.Range("L20").Value = 73.333333333333
$ws.Range("M20").Value = 116.666666666667
$ws.Range("N20").Value = -73.737373737373

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 160
$ws.Range("G21").Value = 199
$ws.Range("H21").Value = -19.597989949748
$ws.Range("I21").Value = 118
$ws.Range("J21").Value = 138
$ws.Range("K21").Value = -14.492753623188
$ws.Range("L21").Value = -2.479338842975
$ws.Range("M21").Value = -11.940298507462
$ws.Range("N21").Value = -73.951434878587

# Row 23
$ws.Range("L23").Value = 0
$ws.Range("L23").NumberFormat = '#,##0.0;"-"#,##0.0'

# Row 24
$ws.Range("C24").Value = 47
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = -16.071428571428
$ws.Range("F24").Value = 176
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = -8.333333333333
$ws.Range("I24").Value = 140
$ws.Range("J24").Value = 157
$ws.Range("K24").Value = -10.828025477707
$ws.Range("L24").Value = 79.487179487179
$ws.Range("M24").Value = 77.215189873417

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -20.408163265306
$ws.Range("I25").Value = 36
$ws.Range("J25").Value = 39
$ws.Range("K25").Value = -7.692307692307
$ws.Range("L25").Value = 38.461538461538
$ws.Range("M25").Value = 0

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = -66.666666666666
$ws.Range("L26").Value = -33.333333333333

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = 100
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 150
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 166.666666666667
$ws.Range("L27").Value = 166.666666666667

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = 0
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 2
$ws.Range("J28").Value = 1
$ws.Range("J28").NumberFormat = '#,##0'
$ws.Range("K28").Value = 100
$ws.Range("K28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L28").Value = 100
$ws.Range("N28").Value = -50

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = '#,##0'
$ws.Range("E29").Value = 0
$ws.Range("E29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 1
$ws.Range("J29").NumberFormat = '#,##0'
$ws.Range("K29").Value = 100
$ws.Range("K29").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L29").Value = 100
$ws.Range("N29").Value = -33.333333333333
